$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, 44559, 68, 2000, 2000, 2000, 667),
  @(3, 44669, 92, 2500, 3000, 2755, 918),
  @(4, 44221, 50, 2500, 2500, 2500, 833),
  @(5, 44260, 60, 3500, 3500, 3500, 1167),
  @(6, 44165, 68, 3000, 3000, 3000, 1000),
  @(7, 44225, 56, 3000, 3000, 3000, 1000),
  @(8, 44537, 88, 2000, 2200, 2091, 697),
  @(9, 44179, 78, 3000, 3000, 3000, 1000),
  @(10, 44574, 50, 3000, 3000, 3000, 1000),
  @(11, 44557, 104, 2000, 2500, 2260, 753),
  @(12, 44967, 110, 3000, 3300, 3136, 1045),
  @(13, 44292, 40, 3000, 3000, 3000, 1000),
  @(14, 45189, 50, 3000, 3000, 3000, 1000),
  @(15, 45092, 90, 3000, 3500, 3278, 1093),
  @(16, 44627, 78, 3500, 3500, 3500, 1167),
  @(17, 44965, 87, 3000, 3000, 3000, 1000),
  @(18, 44390, 50, 3000, 3000, 3000, 1000),
  @(19, 45117, 56, 3000, 3000, 3000, 1000),
  @(20, 44223, 80, 2500, 3000, 2781, 927),
  @(21, 44804, 85, 3000, 3000, 3000, 1000),
  @(22, 44935, 78, 3000, 3000, 3000, 1000),
  @(23, 44224, 67, 3000, 3000, 3000, 1000),
  @(24, 44187, 65, 3000, 3000, 3000, 1000),
  @(25, 44222, 45, 3000, 3000, 3000, 1000),
  @(26, 44756, 104, 2800, 3000, 2904, 968),
  @(27, 44937, 68, 3500, 3500, 3500, 1167),
  @(28, 44291, 45, 3000, 3000, 3000, 1000),
  @(29, 44992, 45, 4000, 4000, 4000, 1333),
  @(30, 45118, 67, 3000, 3000, 3000, 1000),
  @(31, 44389, 81, 2800, 3000, 2889, 963),
  @(32, 44166, 45, 2500, 2500, 2500, 833),
  @(33, 44340, 54, 3000, 3000, 3000, 1000),
  @(34, 44242, 95, 2500, 3000, 2737, 912),
  @(35, 44845, 80, 2500, 2500, 2500, 833),
  @(36, 44536, 125, 2200, 2200, 2200, 733),
  @(37, 44193, 70, 3000, 3000, 3000, 1000),
  @(38, 44243, 45, 3000, 3000, 3000, 1000)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 4).Value = $row[1]
  $ws.Cells.Item($r, 10).Value = $row[2]
  $ws.Cells.Item($r, 11).Value = $row[3]
  $ws.Cells.Item($r, 12).Value = $row[4]
  $ws.Cells.Item($r, 13).Value = $row[5]
  $ws.Cells.Item($r, 16).Value = $row[6]
}
